$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 2499
$ws.Range("J43").Value = 2499
$ws.Range("L43").Value = 2499
$ws.Range("N43").Value = -2637
$ws.Range("H96").Value = 419.0909
$ws.Range("I96").Value = 360.2857
$ws.Range("K96").Value = 1080.8571
$ws.Range("M96").Value = 292.1428999999998
$ws.Range("H100").Value = 4035.1333
$ws.Range("I100").Value = 2393.1667
$ws.Range("J100").Value = 5129.778
$ws.Range("K100").Value = 2393.1667
$ws.Range("L100").Value = 5129.778
$ws.Range("M100").Value = -1852.1667
$ws.Range("N100").Value = -6211.778
$ws.Range("H138").Value = 3295.169
$ws.Range("I138").Value = 3894.842
$ws.Range("J138").Value = 3098.724
$ws.Range("K138").Value = 11684.526
$ws.Range("L138").Value = 9296.172
$ws.Range("M138").Value = -6544.526
$ws.Range("N138").Value = -19576.172

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 172.23077
$ws.Range("I5").Value = 158.09091
$ws.Range("K5").Value = 158.09091
$ws.Range("M5").Value = -46.09091000000001
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H97").Value = 5432.8276
$ws.Range("I97").Value = 5716.6816
$ws.Range("K97").Value = 5716.6816
$ws.Range("M97").Value = -5220.6816
$ws.Range("H109").Value = 74979.2
$ws.Range("J109").Value = 74979.2
$ws.Range("L109").Value = 74979.2
$ws.Range("N109").Value = -77753.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 172.23077
$ws.Range("I4").Value = 158.09091
$ws.Range("K4").Value = 158.09091
$ws.Range("M4").Value = -43.09091000000001
$ws.Range("H99").Value = 2029.4546
$ws.Range("I99").Value = 1993.4
$ws.Range("K99").Value = 1993.4
$ws.Range("M99").Value = -495.4000000000001
$ws.Range("H108").Value = 97784
$ws.Range("J108").Value = 97784
$ws.Range("L108").Value = 97784
$ws.Range("N108").Value = -105464
$ws.Range("H134").Value = 8484.404
$ws.Range("I134").Value = 7667.1035
$ws.Range("K134").Value = 23001.3105
$ws.Range("M134").Value = -20466.3105

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 491.16666
$ws.Range("I22").Value = 449
$ws.Range("J22").Value = 512.25
$ws.Range("K22").Value = 449
$ws.Range("L22").Value = 512.25
$ws.Range("M22").Value = -99
$ws.Range("N22").Value = -1212.25
$ws.Range("H58").Value = 2249.6
$ws.Range("I58").Value = 2417.0908
$ws.Range("J58").Value = 1789
$ws.Range("K58").Value = 2417.0908
$ws.Range("L58").Value = 1789
$ws.Range("M58").Value = -2214.0908
$ws.Range("N58").Value = -2195
$ws.Range("H99").Value = 1511.1177
$ws.Range("I99").Value = 1400.4286
$ws.Range("J99").Value = 2027.6666
$ws.Range("K99").Value = 1400.4286
$ws.Range("L99").Value = 2027.6666
$ws.Range("M99").Value = 97.57140000000004
$ws.Range("N99").Value = -5023.6666
$ws.Range("H105").Value = 1557.3636
$ws.Range("I105").Value = 1194.4286
$ws.Range("K105").Value = 1194.4286
$ws.Range("M105").Value = 552.5714
$ws.Range("H107").Value = 2862.5908
$ws.Range("I107").Value = 3065.2104
$ws.Range("J107").Value = 1579.3334
$ws.Range("K107").Value = 3065.2104
$ws.Range("L107").Value = 1579.3334
$ws.Range("M107").Value = -1145.2104
$ws.Range("N107").Value = -5419.3334
$ws.Range("H126").Value = 1511.1177
$ws.Range("I126").Value = 1400.4286
$ws.Range("J126").Value = 2027.6666
$ws.Range("K126").Value = 4201.2858
$ws.Range("L126").Value = 6082.9998
$ws.Range("M126").Value = -1731.2858
$ws.Range("N126").Value = -11022.9998
$ws.Range("H136").Value = 2249.6
$ws.Range("I136").Value = 2417.0908
$ws.Range("J136").Value = 1789
$ws.Range("K136").Value = 7251.2724
$ws.Range("L136").Value = 5367
$ws.Range("M136").Value = -4701.2724
$ws.Range("N136").Value = -10467

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2862.8125
$ws.Range("I97").Value = 1902.5
$ws.Range("K97").Value = 1902.5
$ws.Range("M97").Value = -1406.5
$ws.Range("H102").Value = 2002.9117
$ws.Range("I102").Value = 1770.7667
$ws.Range("J102").Value = 3744
$ws.Range("K102").Value = 1770.7667
$ws.Range("L102").Value = 3744
$ws.Range("M102").Value = -148.7666999999999
$ws.Range("N102").Value = -6988
$ws.Range("H132").Value = 2546.28
$ws.Range("I132").Value = 2529.8635
$ws.Range("J132").Value = 2666.6667
$ws.Range("K132").Value = 7589.5905
$ws.Range("L132").Value = 8000.000100000001
$ws.Range("M132").Value = -5059.5905
$ws.Range("N132").Value = -13060.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 11249.25
$ws.Range("J14").Value = 12499
$ws.Range("L14").Value = 12499
$ws.Range("N14").Value = -12843
$ws.Range("H22").Value = 1849.6666
$ws.Range("I22").Value = 799
$ws.Range("J22").Value = 2059.8
$ws.Range("K22").Value = 799
$ws.Range("L22").Value = 2059.8
$ws.Range("M22").Value = -504
$ws.Range("N22").Value = -2649.8
$ws.Range("H24").Value = 262374.75
$ws.Range("J24").Value = 17250
$ws.Range("L24").Value = 17250
$ws.Range("N24").Value = -17936
$ws.Range("H27").Value = 1849.6666
$ws.Range("I27").Value = 799
$ws.Range("J27").Value = 2059.8
$ws.Range("K27").Value = 799
$ws.Range("L27").Value = 2059.8
$ws.Range("M27").Value = -692
$ws.Range("N27").Value = -2273.8
$ws.Range("H45").Value = 42355
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 42355
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 42355
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -43169
$ws.Range("H46").Value = 1971.2858
$ws.Range("J46").Value = 2624.75
$ws.Range("L46").Value = 2624.75
$ws.Range("N46").Value = -3000.75
$ws.Range("H55").Value = 1497.4615
$ws.Range("J55").Value = 1913.2222
$ws.Range("L55").Value = 1913.2222
$ws.Range("N55").Value = -2259.2222

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 8318
$ws.Range("I21").Value = 530
$ws.Range("K21").Value = 530
$ws.Range("M21").Value = -295
$ws.Range("H24").Value = 35000
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 35000
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 35000
$ws.Range("M24").ClearContents()
$ws.Range("N24").Value = -35460
$ws.Range("H35").Value = 8318
$ws.Range("I35").Value = 530
$ws.Range("K35").Value = 530
$ws.Range("M35").Value = -240
$ws.Range("H37").Value = 19999.666
$ws.Range("J37").Value = 19999.666
$ws.Range("L37").Value = 19999.666
$ws.Range("N37").Value = -20405.666
$ws.Range("H96").Value = 2979.6667
$ws.Range("I96").Value = 1492.5
$ws.Range("J96").Value = 3723.25
$ws.Range("K96").Value = 1492.5
$ws.Range("L96").Value = 3723.25
$ws.Range("M96").Value = -119.5
$ws.Range("N96").Value = -6469.25
$ws.Range("H100").Value = 2190.9285
$ws.Range("I100").Value = 1214.2858
$ws.Range("K100").Value = 2428.5716
$ws.Range("M100").Value = -1887.5716
